$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.062.00"
$ws.Range("D3").Value = "1.645.70"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "1.874.68"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "1.668.78"
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "26.077.03"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "195.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("E25").Value = "  +4.40%  "
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "143.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("E30").Value = "  +1.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0498"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("E34").Value = "  -2.60%  "
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("E36").Value = "  +0.57%  "
$ws.Range("D37").Value = "1.133.56"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("E38").Value = "  -1.09%  "
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").Value = "1.783.37"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("E46").Value = "  +0.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0526"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.415"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("E51").Value = "  +0.04%  "
